$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "as at" date in the intro paragraph (A2)
$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 02 May 2025"

# 2. Insert three new publication rows into the "28 Jul 2025" week block.
#    Work top-down so each insertion point accounts for the rows already
#    shifted down by the previous insert.

# 2a. New row 27: "Offender Accommodation Outcomes..." (before the existing
#     "Offender management statistics quarterly" row)
$ws.Rows.Item(27).Insert()
$ws.Range("A27").Value = "28 Jul 2025"
$ws.Range("B27").Value = "Offender Accommodation Outcomes, update to March 2025"
$ws.Range("C27").Value = "31 July 2025"
$ws.Range("D27").Value = "provisional"
$ws.Range("E27").Value = 31
$ws.Range("F27").Value = "standard"

# 2b. New row 30: "Offender Employment Outcomes..." (before the existing
#     "Safety in the children and young people secure estate" row)
$ws.Rows.Item(30).Insert()
$ws.Range("A30").Value = "28 Jul 2025"
$ws.Range("B30").Value = "Offender Employment Outcomes, update to March 2025"
$ws.Range("C30").Value = "31 July 2025"
$ws.Range("D30").Value = "provisional"
$ws.Range("E30").Value = 31
$ws.Range("F30").Value = "standard"

# 2c. New row 32: "Community Performance Annual..." (before the existing
#     "HMPPS Annual Digest" row)
$ws.Rows.Item(32).Insert()
$ws.Range("A32").Value = "28 Jul 2025"
$ws.Range("B32").Value = "Community Performance Annual, update to March 2025"
$ws.Range("C32").Value = "31 July 2025"
$ws.Range("D32").Value = "provisional"
$ws.Range("E32").Value = 31
$ws.Range("F32").Value = "standard"

# 3. Extend the conditional-formatting ranges to cover the new rows
#    (table now runs through row 61 instead of row 58).
$fcs = $ws.Range("A5:F61").FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws.Range("A5:F61"))
$fcs.Item(2).ModifyAppliesToRange($ws.Range("A5:F61"))
$fcs.Item(3).ModifyAppliesToRange($ws.Range("A5:F61"))
$fcs.Item(4).ModifyAppliesToRange($ws.Range("A5:A61"))
$fcs.Item(5).ModifyAppliesToRange($ws.Range("A5:A61"))
